$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "BlackBox"

# --- Apply a box border ("All Borders") to the data rows of both tables ---
# STAY CLASS data table (rows 3-12)
$ws.Range("A3:D12").Borders.LineStyle = 1
$ws.Range("A3:D12").Borders.Weight = 2

# GATE CLASS data table (rows 19-27, before the row-delete shift below)
$ws.Range("A19:D27").Borders.LineStyle = 1
$ws.Range("A19:D27").Borders.Weight = 2

# --- Highlight the two section title cells with a green fill ---
$ws.Range("A1").Interior.Color = 5296274
$ws.Range("A17").Interior.Color = 5296274

# --- Tidy up leftover placeholder-styled blank cells ---
$ws.Range("E19:F19").Clear()

# --- Shrink the title row height ---
$ws.Rows(1).RowHeight = 18

# --- Remove the extra blank rows between the two tables ---
$ws.Rows("13:14").Delete()

# --- Update the remembered selection ---
$ws.Range("C28").Select()
